# Release v2.3 deployment: the CSV-mapping "conceptPath" column (J) on
# Sheet1 pointed its tree paths at a generic "/root/..." namespace. For
# this dataset the paths must live under "/dementia/..." instead, so every
# occurrence of the "/root/" prefix is swapped for "/dementia/" across the
# whole sheet (28 distinct values, 34 cells once shared duplicates count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPart = 2 (match anywhere within the cell, not the whole cell contents)
$ws.Cells.Replace("/root/", "/dementia/", 2) | Out-Null

# Re-point the on-screen selection/scroll at the edited block, matching
# where the editor was working when the change was made.
$ws.Range("J2,J4:J6,J8:J9,J11:J38,O15").Select() | Out-Null
$ws.Range("O15").Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
